# "Add files via upload" -- populate the tail of the "test" sheet
# (contact / donation / attachment-description strings) and move the
# selection down to where the new content now lives.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

$ws.Range("A39").Value = "Контакт для предложений и пожеланий в развитие бота."
$ws.Range("A40").Value = "📥 Если вы хотите поддержать проект, мы с благодарностью примем Вашу поддержку в любом объеме`n13jtTtFix1ji1j8dzk3WAeo6B1A3hY9FKX"
$ws.Range("A41").Value = "Напишите описание к вложению"
$ws.Range("A42").Value = "Описание добавлено"

$ws.Rows.Item(39).RowHeight = 30
$ws.Rows.Item(40).RowHeight = 105
$ws.Rows.Item(41).RowHeight = 30

$ws.Range("A43").Select()
